$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format before writing, so values like
# "206.83" are stored as text (matching the source inline-string cells)
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.675.55"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.584.06"
$ws.Range("E3").Value = "  -2.51%  "
$ws.Range("E4").Value = "  +0.79%  "
$ws.Range("D5").Value = "206.83"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("E6").Value = "  -3.08%  "
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("D8").Value = "22.13"
$ws.Range("E8").Value = "  -4.67%  "
$ws.Range("D9").Value = "0.251"
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("D10").Value = "0.0590"
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("D11").Value = "0.0865"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").Value = "1.808.14"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").Value = "1.575.37"
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("E14").Value = "  -4.03%  "
$ws.Range("D15").Value = "0.528"
$ws.Range("E15").Value = "  -5.00%  "
$ws.Range("D16").Value = "63.45"
$ws.Range("E16").Value = "  -2.24%  "
$ws.Range("D17").Value = "27.627.58"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "219.17"
$ws.Range("E18").Value = "  -3.84%  "
$ws.Range("D19").Value = "0.0₃0694"
$ws.Range("E19").Value = "  -3.18%  "
$ws.Range("D20").Value = "7.31"
$ws.Range("E20").Value = "  -3.43%  "
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D23").Value = "9.52"
$ws.Range("E23").Value = "  -3.92%  "
$ws.Range("D24").Value = "1.97"
$ws.Range("E24").Value = "  -3.64%  "
$ws.Range("D25").Value = "153.84"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "6.85"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").Value = "15.10"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("E29").Value = "  -4.44%  "
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("D31").Value = "0.0466"
$ws.Range("E31").Value = "  -2.91%  "
$ws.Range("E32").Value = "  -5.82%  "
$ws.Range("D33").Value = "1.360.01"
$ws.Range("E33").Value = "  -2.98%  "
$ws.Range("E34").Value = "  -4.83%  "
$ws.Range("D35").Value = "1.53"
$ws.Range("E35").Value = "  -4.47%  "
$ws.Range("D36").Value = "0.966"
$ws.Range("E36").Value = "  -3.10%  "
$ws.Range("D37").Value = "2.31"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").Value = "0.0167"
$ws.Range("E38").Value = "  -1.43%  "
$ws.Range("D39").Value = "0.534"
$ws.Range("E39").Value = "  -3.47%  "
$ws.Range("D40").Value = "0.819"
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("D42").Value = "0.971"
$ws.Range("E42").Value = "  -2.70%  "
$ws.Range("D43").Value = "63.67"
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("E44").Value = "  +2.91%  "
$ws.Range("E45").Value = "  -4.01%  "
$ws.Range("E46").Value = "  -4.08%  "
$ws.Range("D47").Value = "1.720.23"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("D48").Value = "88.25"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").Value = "0.0₆0100"
$ws.Range("E49").Value = "  +12.11%  "
$ws.Range("D50").Value = "0.0969"
$ws.Range("E50").Value = "  -3.74%  "
$ws.Range("D51").Value = "0.0497"
$ws.Range("E51").Value = "  -1.13%  "

# Restore the default cell style on column D (removes the Text
# number-format override so cells keep their original "Normal" style).
$ws.Range("D2:D51").Style = "Normal"
